# Generate Report for Handback
#
# The localization-status workbook is updated to reflect that the handback
# xliff files have come back "in sync" with en-US for both the zh-cn and
# de-de locales:
#   - Overview sheet: Status columns (zh-cn / de-de) flip from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - zh-cn / de-de sheets: Status flips the same way, the "Latest Target
#     File" + "Latest Handback File" cells get populated (with the target
#     file becoming a hyperlink back to the source .md, like column A), and
#     "Latest Handback DateTime" gets stamped with the real handback time.

$wb = $excel.ActiveWorkbook

$sourceMdName = "8065581c-559f-45a9-a175-93f3cbafd4ca.md"
$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2e22ab09f38dc91a28a4b83fced5f94311f9f2e/e2e/8065581c-559f-45a9-a175-93f3cbafd4ca.md"

$statusInSync = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusInSync
$wsOverview.Range("F2").Value = $statusInSync

# Status column auto-fit wider to hold the new, longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.125
$wsOverview.Columns.Item(6).ColumnWidth = 29.125

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusInSync
$wsZhCn.Range("I2").Value = $sourceMdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $sourceUrl, "", "", $sourceMdName) | Out-Null
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("J2").Value = "8065581c-559f-45a9-a175-93f3cbafd4ca.829fbc2bb48ea4a57977d7fccd56083cce33b0ff.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-21 01:03:55"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.125
$wsZhCn.Columns.Item(9).ColumnWidth = 39.125
$wsZhCn.Columns.Item(10).ColumnWidth = 39.125

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusInSync
$wsDeDe.Range("I2").Value = $sourceMdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $sourceUrl, "", "", $sourceMdName) | Out-Null
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("J2").Value = "8065581c-559f-45a9-a175-93f3cbafd4ca.829fbc2bb48ea4a57977d7fccd56083cce33b0ff.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-21 01:04:04"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.125
$wsDeDe.Columns.Item(9).ColumnWidth = 39.125
$wsDeDe.Columns.Item(10).ColumnWidth = 39.125
